$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log row (row 7) with the Testmail #4 data
$ws.Cells.Item(7, 1).Value = "Wil je dit artikel voor me inkopen?"
$ws.Cells.Item(7, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(7, 3).Value = "Testmail #4: Wil je dit artikel voor me inkopen?"
$ws.Cells.Item(7, 4).Value = "Productinformatie"
$ws.Cells.Item(7, 5).Value = "Beste klant,`nBedankt voor je interesse in ons artikel. Helaas kan ik je op basis van dit bericht niet verder helpen. Kun je meer details geven over welk artikel je wilt inkopen en op welke manier? Zo kan ik je beter assisteren.`nMet vriendelijke groet,`n[Jouw naam]`nE-mailassistent"
$ws.Cells.Item(7, 6).Value = "2025-06-29 14:06:13"
$ws.Cells.Item(7, 7).Value = "Ja"
$ws.Cells.Item(7, 8).Value = "Ja"
$ws.Cells.Item(7, 9).Value = "Nee"

# Undo the automatic row-height/auto-fit side effect triggered by the
# multi-line value so the row keeps its default (non-custom) height.
$ws.Rows.Item(7).AutoFit()

# Extend the conditional formatting ranges on the Logs sheet to include row 7
$dFcs = $ws.Range("D2:D6").FormatConditions
for ($i = 1; $i -le $dFcs.Count; $i++) {
    $dFcs.Item($i).ModifyAppliesToRange($ws.Range("D2:D7"))
}

$gFcs = $ws.Range("G2:G6").FormatConditions
for ($i = 1; $i -le $gFcs.Count; $i++) {
    $gFcs.Item($i).ModifyAppliesToRange($ws.Range("G2:G7"))
}

$hFcs = $ws.Range("H2:H6").FormatConditions
for ($i = 1; $i -le $hFcs.Count; $i++) {
    $hFcs.Item($i).ModifyAppliesToRange($ws.Range("H2:H7"))
}

$iFcs = $ws.Range("I2:I6").FormatConditions
for ($i = 1; $i -le $iFcs.Count; $i++) {
    $iFcs.Item($i).ModifyAppliesToRange($ws.Range("I2:I7"))
}

# Update the Dashboard sheet: Productinformatie count goes from 1 to 2
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(4, 2).Value = 2
